# "contingencies with rene fine"
# Two new line rows (line7 / line8) were inserted into the underlying
# shared-string table right before the "extr*" block. That insertion
# shifts every row below it down one "slot" in the B column, which is
# why rows 8-17 all end up pointing at a different label than before,
# and two brand new rows (16 and 17) appear at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B relabeling caused by the shared-string insertion ---
$ws.Cells.Item(8, 2).Value  = "line7"
$ws.Cells.Item(9, 2).Value  = "line8"
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(15, 2).Value = "extr6"

# --- Updated numeric / boolean data for existing rows 8-15 ---
$ws.Cells.Item(8, 3).Value  = 14
$ws.Cells.Item(8, 4).Value  = 11
$ws.Cells.Item(8, 5).Value  = $true

$ws.Cells.Item(9, 3).Value  = 16
$ws.Cells.Item(9, 4).Value  = 9
$ws.Cells.Item(9, 5).Value  = $true

$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# --- Two brand new rows appended at the bottom ---
# Clone column A's formatting (bold / centered / bordered style) from the
# row above so the new index cells match the rest of the table.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A15").Copy($ws.Range("A17"))

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $true

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
